# session_1 intro slides: update module timing estimates on the
# "Let's get started" agenda slide.
#   Module 4 (Basics of Visual Studio Code Git integration): (15 min)  -> (5-10 min)
#   Module 5 (What makes a good repository):                  (30 min) -> (20-30 min)

$p = $ppt.ActivePresentation

# Locate the agenda slide/shape robustly (don't hard-code indices): the
# shape whose text contains the module list.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -like "*Module 4*" -and $t -like "*Module 5*") {
                $targetSlide = $s
                $targetShape = $sh
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# --- Module 4 timing: (15 min) -> (5-10 min) ---
$text = $tr.Text
$idx15 = $text.IndexOf("(15 min)")
if ($idx15 -ge 0) {
    $run15 = $tr.Characters($idx15 + 1, 8)
    $run15.Text = "(5-10 min)"
}

# --- Module 5 timing: the (30 min) occurring after (15 min) -> (20-30 min) ---
$text = $tr.Text
$idx30 = $text.IndexOf("(30 min)", $idx15)
if ($idx30 -ge 0) {
    $run30 = $tr.Characters($idx30 + 1, 8)
    $run30.Text = "(20-30 min)"
}

Write-Host "Updated agenda text:"
Write-Host $tr.Text
